$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.158.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.05%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +3.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.62%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +3.62%  "

# Row 10 - OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.89%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.26%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.21%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.355.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.66%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "

# Row 15 - Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.27%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +2.06%  "

# Row 17 - Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "

# Row 18 - WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.053.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "

# Row 19 - WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.108.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.78%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.49%  "

# Row 21 - Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.55%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +1.77%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "224.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.22%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.06%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +1.82%  "

# Row 26 - PancakeSwap
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.17%  "

# Row 28 - Cosmos->ImmutableX
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.15%  "

# Row 29 - ImmutableX->Cosmos
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "

# Row 30 - EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

# Row 31 - Kaspa
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.125"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "

# Row 32 - Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.59%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +2.50%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.21%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +1.25%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  +0.11%  "

# Row 38 - RenderToken->WEMIXToken
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.81%  "

# Row 39 - WEMIXToken->RenderToken
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.18%  "

# Row 40 - THORChain
$ws.Range("E40").Value = "  -1.05%  "

# Row 41 - HuobiToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

# Row 42 - FTXToken->Maker
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.483.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "

# Row 43 - Maker->FTXToken
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.88%  "

# Row 44 - Aave
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.52%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.04%  "

# Row 46 - Cronos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0929"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +3.20%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +1.61%  "

# Row 49 - InjectiveProtocol->FraxShare
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.71%  "

# Row 50 - FraxShare->InjectiveProtocol
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "15.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.85%  "

# Row 51 - MXToken
$ws.Range("E51").Value = "  +1.70%  "
